# chore: adapt column header formatting to respective input file names
#
# - Rename the diff-table header columns from the generic "_old" / "_new"
#   suffixes to the concrete format-version suffixes "_FV2404" / "_FV2410".
# - Freeze the header row (row 1).
# - Wrap the data range in an Excel Table ("Table1") with an AutoFilter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename header row (A1:U1): "<Name>_old" -> "<Name>_FV2404"
#    and "<Name>_new" -> "<Name>_FV2410". Column K ("diff") is untouched.
# ---------------------------------------------------------------------
$oldHeaders = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old"
)

$newHeaders = @(
    "Segmentname_new",
    "Segmentgruppe_new",
    "Segment_new",
    "Datenelement_new",
    "Segment ID_new",
    "Code_new",
    "Qualifier_new",
    "Beschreibung_new",
    "Bedingungsausdruck_new",
    "Bedingung_new"
)

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i] -replace "_old$", "_FV2404"
}

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i] -replace "_new$", "_FV2410"
}

# ---------------------------------------------------------------------
# 2) Freeze the header row.
# ---------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------
# 3) Turn the used range into an Excel Table ("Table1") with a header row.
# ---------------------------------------------------------------------
$dataRange = $ws.Range("A1:U89")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1, $null)
$tbl.Name = "Table1"
